# cargaUsuarios.xlsx - corrige el tipo de documento de Maria Garcia (fila 3)
# de "CE" a "CC", y deja seleccionada toda la tabla de datos (A1:H4) tal
# como quedo el libro al guardarlo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 3 = Maria Garcia -> columna F ("tipo documento"): CE -> CC
$ws.Range("F3").Value = "CC"

# Selecciona la tabla completa de datos (A1:H4)
$ws.Range("A1:H4").Select()
